$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 2 (shifting existing data rows down by 2)
$ws.Rows.Item(2).Resize(2).Insert()

# The freshly inserted rows pick up formatting (bold) from the row above (the
# header); clear that and reapply the same formatting used by the rest of the
# data rows (plain cells, with column A formatted as a date/time value).
$ws.Range("A2:I3").ClearFormats()
$ws.Range("A2:A3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the two new rows with the latest lottery draw
$ws.Range("A2").Value = 45356
$ws.Range("B2").Value = "03,18,24,35,40,41"
$ws.Range("C2").Value = "Within"
$ws.Range("D2").Value = "Within"
$ws.Range("E2").Value = "Within"
$ws.Range("F2").Value = "Within"
$ws.Range("G2").Value = "Within"
$ws.Range("H2").Value = "Within"
$ws.Range("I2").Value = "day_name"

$ws.Range("A3").Value = 45356
$ws.Range("B3").Value = "03,18,24,35,40,41"
$ws.Range("C3").Value = "Within"
$ws.Range("D3").Value = "Within"
$ws.Range("E3").Value = "Within"
$ws.Range("F3").Value = "Within"
$ws.Range("G3").Value = "Within"
$ws.Range("H3").Value = "Within"
$ws.Range("I3").Value = "month_day"

$wb.Save()
